# Update "想去人数" (F column) figures across the sheets, as per the
# site regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 8751
$ws1.Range("F17").Value = 1406
$ws1.Range("F21").Value = 1328

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 13
$ws2.Range("F27").Value = 169

# --- Sheet: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 260
$ws3.Range("F8").Value = 1959

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 8751
$ws4.Range("F8").Value = 260
$ws4.Range("F11").Value = 1959
$ws4.Range("F13").Value = 13
$ws4.Range("F23").Value = 1406
$ws4.Range("F25").Value = 1328
$ws4.Range("F38").Value = 169
